$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated values (new TPM-based computation) for columns E:T, rows 2-4
$data = @{
    2 = @(1, 0.3333333333333333, 0.05881766666666666, 0.176453, 1, 1, 3, 1, 2.473670666666667, 7.421011999999999, 0.04625390589638428, 0.04625390589638428, 0.1454955367151111, 1.309459830436, 0.04625390589638428, 0.04625390589638428)
    3 = @(1, 0.3333333333333333, 0.05881766666666666, 0.176453, 1, 1, 3, 1, 30.598211, 91.794633, 0.5721403383494235, 0.5721403383494235, 1.799715375194333, 16.197438376749, 0.5721403383494235, 0.5721403383494235)
    4 = @(1, 0.3333333333333333, 0.05881766666666666, 0.176453, 1, 1, 3, 1, 20.408373, 61.225119, 0.3816057557541923, 0.3816057557541923, 1.200372880323, 10.803355922907, 0.3816057557541923, 0.3816057557541923)
}

# Columns E (5) through T (20)
$startCol = 5

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $values.Length; $i++) {
        $col = $startCol + $i
        $ws.Cells.Item($row, $col).Value = $values[$i]
    }
}
